$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'63.092.65"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.77%  "

$c = $ws.Range("D3")
$c.Value = "'3.173.38"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -4.36%  "

$ws.Range("E4").Value = "  +0.10%  "

$c = $ws.Range("D5")
$c.Value = "'591.01"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.38%  "

$c = $ws.Range("D6")
$c.Value = "'136.26"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -3.99%  "

$ws.Range("E7").Value = "  -0.01%  "

$c = $ws.Range("D8")
$c.Value = "'3.169.24"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -4.39%  "

$c = $ws.Range("D9")
$c.Value = "'0.512"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.39%  "

$ws.Range("E10").Value = "  -4.53%  "

$ws.Range("E11").Value = "  -3.60%  "

$c = $ws.Range("D12")
$c.Value = "'0.456"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -2.51%  "

$c = $ws.Range("D14")
$c.Value = "'34.88"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.55%  "

$c = $ws.Range("D15")
$c.Value = "'3.696.76"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -4.34%  "

$ws.Range("E16").Value = "  -2.03%  "

$c = $ws.Range("D17")
$c.Value = "'3.173.06"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -4.39%  "

$c = $ws.Range("D18")
$c.Value = "'63.024.59"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.01%  "

$c = $ws.Range("D19")
$c.Value = "'6.62"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -3.49%  "

$c = $ws.Range("D20")
$c.Value = "'461.36"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -3.84%  "

$c = $ws.Range("D21")
$c.Value = "'13.93"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.44%  "

$ws.Range("E22").Value = "  -3.37%  "

$c = $ws.Range("D23")
$c.Value = "'7.66"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -6.51%  "

$c = $ws.Range("D24")
$c.Value = "'13.46"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.80%  "

$c = $ws.Range("D25")
$c.Value = "'83.40"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.82%  "

$c = $ws.Range("D26")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("E27").Value = "  -2.96%  "

$ws.Range("E28").Value = "  +0.12%  "

$ws.Range("E29").Value = "  -4.34%  "

$c = $ws.Range("D30")
$c.Value = "'6.81"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -5.83%  "

$ws.Range("E31").Value = "  -5.93%  "

$c = $ws.Range("D32")
$c.Value = "'27.34"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -5.51%  "

$c = $ws.Range("D33")
$c.Value = "'0.103"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.84%  "

$ws.Range("E35").Value = "  -6.05%  "

$c = $ws.Range("D36")
$c.Value = "'5.83"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -3.63%  "

$c = $ws.Range("D37")
$c.Value = "'51.16"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.88%  "

$ws.Range("E38").Value = "  -4.07%  "

$ws.Range("E39").Value = "  -2.34%  "

$c = $ws.Range("D40")
$c.Value = "'405.81"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -6.51%  "

$c = $ws.Range("D41")
$c.Value = "'2.68"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.40%  "

$c = $ws.Range("D42")
$c.Value = "'8.11"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.75%  "

$ws.Range("E43").Value = "  -3.59%  "

$c = $ws.Range("D44")
$c.Value = "'2.787.49"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -9.84%  "

$ws.Range("E45").Value = "  -3.79%  "

$ws.Range("E46").Value = "  -2.40%  "

$c = $ws.Range("D48")
$c.Value = "'25.72"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.21%  "

$c = $ws.Range("D49")
$c.Value = "'34.64"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -5.90%  "

$c = $ws.Range("D50")
$c.Value = "'122.67"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.73%  "

$ws.Range("E51").Value = "  -2.12%  "
